$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Femacal de La Calera" / Zapallo
# italiano. It is inserted as row 742 (pushing the existing rows 742:770
# down to 743:771, dimension A1:R770 -> A1:R771).
$ws.Rows.Item(742).Insert()

$ws.Cells.Item(742, 1).Value  = 3
$ws.Cells.Item(742, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(742, 3).Value  = "Coquimbo"
$ws.Cells.Item(742, 4).Value  = 45075
$ws.Cells.Item(742, 5).Value  = 5
$ws.Cells.Item(742, 6).Value  = 100112032
$ws.Cells.Item(742, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(742, 8).Value  = "Sin especificar"
$ws.Cells.Item(742, 9).Value  = "Primera"
$ws.Cells.Item(742, 10).Value = 210
$ws.Cells.Item(742, 11).Value = 8500
$ws.Cells.Item(742, 12).Value = 9000
$ws.Cells.Item(742, 13).Value = 8755
$ws.Cells.Item(742, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(742, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(742, 16).Value = 146
$ws.Cells.Item(742, 17).Value = 60
$ws.Cells.Item(742, 18).Value = "Hortaliza"
